# Common: Initial wizard UI stuff
# Update the "Import" sheet (translations fixture) with the new
# lab.wizard.build.* keys, replacing the old lab.wizard.build.title /
# lab.wizard.build.subtitle rows and appending the new wizard-tab rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Import")

# --- Update existing rows 189-190 (previously title/subtitle) ---
$ws.Range("B189").Value = "lab.wizard.build.first.title"
$ws.Range("C189").Value = "Průvodce novým buildem"

$ws.Range("B190").Value = "lab.wizard.build.first.subtitle"
$ws.Range("C190").Value = "Tento průvodce vám pomůže zaevidovat nový build."

# --- Append new rows 191-200, reusing row 190's formatting (style s="1") ---
$ws.Range("A190:C190").Copy()
$ws.Range("A191:C200").PasteSpecial(-4122)

$ws.Range("A191").Value = "cs"
$ws.Range("B191").Value = "lab.wizard.build.first.tab.label"
$ws.Range("C191").Value = "Úvod"

$ws.Range("A192").Value = "cs"
$ws.Range("B192").Value = "lab.wizard.build.first.tab.description"
$ws.Range("C192").Value = "Průvodce novým buildem"

$ws.Range("A193").Value = "cs"
$ws.Range("B193").Value = "lab.wizard.build.coil.tab.label"
$ws.Range("C193").Value = "Spirálka"

$ws.Range("A194").Value = "cs"
$ws.Range("B194").Value = "lab.wizard.build.coil.tab.description"
$ws.Range("C194").Value = "Vyberte prosím použitou spirálku"

$ws.Range("A195").Value = "cs"
$ws.Range("B195").Value = "lab.wizard.build.atomizer.tab.label"
$ws.Range("C195").Value = "Atomizér"

$ws.Range("A196").Value = "cs"
$ws.Range("B196").Value = "lab.wizard.build.atomizer.tab.description"
$ws.Range("C196").Value = "Vyberte prosím použitý atomizér"

$ws.Range("A197").Value = "cs"
$ws.Range("B197").Value = "lab.wizard.build.cotton.tab.label"
$ws.Range("C197").Value = "Vata"

$ws.Range("A198").Value = "cs"
$ws.Range("B198").Value = "lab.wizard.build.cotton.tab.description"
$ws.Range("C198").Value = "Vyberte prosím vatu"

$ws.Range("A199").Value = "cs"
$ws.Range("B199").Value = "lab.wizard.build.build.tab.label"
$ws.Range("C199").Value = "Build"

$ws.Range("A200").Value = "cs"
$ws.Range("B200").Value = "lab.wizard.build.build.tab.description"
$ws.Range("C200").Value = "Doplňující informace o buildu"

# --- Match the reported viewport / selection state ---
$ws.Range("B193").Select()
